$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values: B/C/D/E take on the values that were in O/R/AN/AQ
$ws.Range("B1").Value2 = $ws.Range("O1").Value2
$ws.Range("C1").Value2 = $ws.Range("R1").Value2
$ws.Range("D1").Value2 = $ws.Range("AN1").Value2
$ws.Range("E1").Value2 = $ws.Range("AQ1").Value2

# Row 2 values: B/C/D/E take on the values that were in O/R/AN/AQ
$ws.Range("B2").Value2 = $ws.Range("O2").Value2
$ws.Range("C2").Value2 = $ws.Range("R2").Value2
$ws.Range("D2").Value2 = $ws.Range("AN2").Value2
$ws.Range("E2").Value2 = $ws.Range("AQ2").Value2

# Row 3 values: B/C/D/E take on the values that were in O/R/AN/AQ
$ws.Range("B3").Value2 = $ws.Range("O3").Value2
$ws.Range("C3").Value2 = $ws.Range("R3").Value2
$ws.Range("D3").Value2 = $ws.Range("AN3").Value2
$ws.Range("E3").Value2 = $ws.Range("AQ3").Value2

# Update selection to match the narrowed range of interest
[void]$ws.Range("B1:E3").Select()
